# Generate Report for Handoff
# - Updates "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#   for the 223f1c98-fa34-48ac-bf8e-dbf83b875c17.md row set (rows 8,9,11,12,13,14)
# - Sets the "Priority" column to "ht" for those same rows on the zh-cn and de-de sheets

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 11, 12, 13, 14)

# --- Overview sheet: column G holds "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-24 02:21:33"
}

# --- zh-cn sheet: column H holds "Latest Handoff Datetime", column E holds "Priority" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-24 02:21:28"
    $wsZhCn.Range("E$r").Value = "ht"
}

# --- de-de sheet: column H holds "Latest Handoff Datetime", column E holds "Priority" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-08-24 02:21:33"
    $wsDeDe.Range("E$r").Value = "ht"
}
